# Applies the "Updated symbol list" crypto-price refresh (Sun Dec 18 2022 14:08 UTC run).
# Source cells are stored as text (inlineStr) even though most look numeric, so
# numeric-looking updates are written with a leading quote (forces text entry, like
# typing '245.81 into Excel) and then restyled back to Normal so no stray text-format
# style sticks to the cell -- only the cell text itself should change, matching the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updates for cells whose new text looks numeric (Price column D, Hora column G) ---
$numericUpdates = [ordered]@{
    'D2' = '245.81'
    'G2' = '14'
    'D3' = '22.41'
    'G3' = '14'
    'D4' = '5.538'
    'G4' = '14'
    'D5' = '0.05633'
    'G5' = '14'
    'D6' = '6.473'
    'G6' = '14'
    'D7' = '0.8047'
    'G7' = '14'
    'D8' = '1.058'
    'G8' = '14'
    'D9' = '0.1464'
    'G9' = '14'
    'D10' = '0.07341'
    'G10' = '14'
    'D11' = '0.03196'
    'G11' = '14'
    'D12' = '0.02928'
    'G12' = '14'
    'D13' = '0.09257'
    'G13' = '14'
    'D14' = '0.001674'
    'G14' = '14'
    'D15' = '3.198'
    'G15' = '14'
    'D16' = '0.04729'
    'G16' = '14'
    'D17' = '0.0005842'
    'G17' = '14'
    'D18' = '0.006279'
    'G18' = '14'
    'D19' = '0.001055'
    'G19' = '14'
    'D20' = '0.004110'
    'G20' = '14'
    'D21' = '0.0001503'
    'G21' = '14'
    'D22' = '3.970'
    'G22' = '14'
    'D23' = '3.383'
    'G23' = '14'
    'D24' = '2.140'
    'G24' = '14'
    'G25' = '14'
    'D26' = '0.1317'
    'G26' = '14'
    'G27' = '14'
    'G28' = '14'
    'G29' = '14'
    'G30' = '14'
    'G31' = '14'
    'G32' = '14'
    'G33' = '14'
    'G34' = '14'
    'G35' = '14'
    'G36' = '14'
    'G37' = '14'
    'G38' = '14'
    'G39' = '14'
    'D40' = '0.04162'
    'G40' = '14'
    'D41' = '0.006885'
    'G41' = '14'
    'D42' = '0.003508'
    'G42' = '14'
    'D43' = '0.1039'
    'G43' = '14'
    'D44' = '0.009842'
    'G44' = '14'
    'D45' = '0.00005658'
    'G45' = '14'
    'G46' = '14'
    'G47' = '14'
    'D48' = '0.02161'
    'G48' = '14'
    'D49' = '0.00002105'
    'G49' = '14'
    'G50' = '14'
    'G51' = '14'
}

foreach ($cellRef in $numericUpdates.Keys) {
    $range = $ws.Range($cellRef)
    $range.Value = "'" + $numericUpdates[$cellRef]
    $range.Style = "Normal"
}

# --- Updates for cells whose new text is plain text (Coin, Link, Volume label columns) ---
$textUpdates = [ordered]@{
    'E17' = '16OneONEWorstin24h'
    'B41' = 'KickToken'
    'C41' = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
    'E41' = '40KickTokenKICK'
    'B43' = 'BKEXToken'
    'C43' = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
    'E43' = '42BKEXTokenBKK'
}

foreach ($cellRef in $textUpdates.Keys) {
    $ws.Range($cellRef).Value = $textUpdates[$cellRef]
}

